# edit.ps1
# Applies the commit's change to orcamentos.xlsx:
#   - "quotations" sheet: a new quotation (BANCO SANTANDER / visita tecnica) is
#     inserted as row 15 (rows 15-28 shift down to 16-29), and the old last
#     row (the "TESTE" quotation) is removed.
#   - "items" sheet: the matching new line-item ("Hora de trabalho ECO") is
#     inserted as row 62 (rows 62-101 shift down to 63-102), and the old last
#     row (the "REFRIGERADOR HORIZONTAL TOPEMA" item, which belonged to the
#     removed "TESTE" quotation) is removed.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    # Force the value to be stored as TEXT even when it looks numeric /
    # date-like (matches the source export, which always serialises these
    # fields as strings), then strip the temporary "@" format so the cell
    # keeps the sheet's default (unstyled) look.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "quotations"
# ---------------------------------------------------------------------
$qs = $wb.Worksheets.Item("quotations")

# Insert the new row at position 15 (pushes old rows 15..29 to 16..30).
$qs.Rows.Item(15).Insert()

Set-TextValue $qs.Range("A15") "NWQzZGY0ZDgtNDMyNS00ZDhlLWEwNDUtOGNlYjkxNzRmM2Y3OjU3MDE2"
Set-TextValue $qs.Range("B15") "RKDOWJG7-E"
Set-TextValue $qs.Range("C15") "90400888000142 - BANCO SANTANDER (BRASIL) S/A"
Set-TextValue $qs.Range("D15") "REFERENTE A VISITA TECNICA PARA DIAGNOSTICO `nO VALOR PODE SER ALTERADO CASO PRECISE DE MAIS TEMPO PARA MANUNTEN€AO`nPE€AS NAO ESTAO INCLUSAS NESTE OR€AMENTO , CASO SEJA NECESSARIO, SERA ADICIONADO NO VALOR FINAL".Replace("€", "Ç")
$qs.Range("F15").Value = $false
Set-TextValue $qs.Range("G15") "350"
Set-TextValue $qs.Range("H15") "350"
Set-TextValue $qs.Range("I15") "Pendente"
Set-TextValue $qs.Range("J15") "2025-12-24T21:32:57.350Z"
Set-TextValue $qs.Range("M15") "Adriana Vieira Masini"
Set-TextValue $qs.Range("O15") "2025-12-17T21:33:20.182Z"
Set-TextValue $qs.Range("P15") "ODAxNDQ5MTMtNDYwMi00MjhmLWE1MWUtMWY1Y2I5NGIxY2Y1OjU3MDE2"
Set-TextValue $qs.Range("Q15") "percentage"
Set-TextValue $qs.Range("R15") "0"
Set-TextValue $qs.Range("S15") "0"
Set-TextValue $qs.Range("T15") "MzMzMDQ1Mjo1NzAxNg=="
Set-TextValue $qs.Range("U15") "pending"

# Remove the old trailing "TESTE" quotation row, now shifted to row 30.
$qs.Rows.Item(30).Delete()

# ---------------------------------------------------------------------
# Sheet "items"
# ---------------------------------------------------------------------
$its = $wb.Worksheets.Item("items")

# Insert the matching new line item at position 62 (pushes old rows
# 62..101 to 63..102).
$its.Rows.Item(62).Insert()

Set-TextValue $its.Range("A62") "MmUwOWMzZjUtNzljZS00OTZhLWExMTQtMDNhNTZlMmMxZDViOjU3MDE2"
$its.Range("B62").Value = 1
$its.Range("C62").Value = 350
Set-TextValue $its.Range("D62") "Hora de trabalho ECO"
$its.Range("E62").Value = 2
Set-TextValue $its.Range("F62") "NWQzZGY0ZDgtNDMyNS00ZDhlLWEwNDUtOGNlYjkxNzRmM2Y3OjU3MDE2"
Set-TextValue $its.Range("G62") "ODY3OTE5NTMtMDdjZi00YzM1LThkN2QtNDc5NzNmNzVkMGY0OjU3MDE2"
$its.Range("H62").Value = 350
Set-TextValue $its.Range("I62") "service"
Set-TextValue $its.Range("J62") "NWQzZGY0ZDgtNDMyNS00ZDhlLWEwNDUtOGNlYjkxNzRmM2Y3OjU3MDE2"

# Remove the old trailing "REFRIGERADOR HORIZONTAL TOPEMA" item row
# (belonged to the removed "TESTE" quotation), now shifted to row 102.
$its.Rows.Item(102).Delete()
